$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need the cell
# format forced to Text first, otherwise Excel auto-converts the
# literal (e.g. "97.00" -> 97, "587.29" -> 587.28999999999996) and
# the text formatting the site displays (trailing zeros, exact
# decimal digits) would be lost.

$ws.Range("D2").Value = "68.959.18"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "3.518.73"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.29"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.72"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").Value = "3.509.78"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.189"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.87"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.56"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "4.077.64"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "626.47"
$ws.Range("E16").Value = "  -6.27%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.49"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "3.517.53"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "68.991.25"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.16"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.888"
$ws.Range("E23").Value = "  -4.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.92"
$ws.Range("E24").Value = "  -6.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.00"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.83"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -4.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -6.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.72"
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("E34").Value = "  -5.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "639.95"
$ws.Range("E35").Value = "  +11.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.74"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  -11.32%  "
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.23"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0454"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("D43").Value = "3.389.22"
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.80"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.41"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("E51").Value = "  +14.88%  "
